$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM data for rows 2-10 (columns A-T)
# Row 8 previously held the last existing row (Target=Resolving-Mac); the data
# is now fully reshuffled: "Resolving-Mac" becomes a new Sending cluster (rows 8-10)
# and is no longer used as a Target cluster. All numeric TPM-derived values were recomputed.
$rows = @(
    @("FAPs", "Fgf17", "Fgfr2", "ECs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 2.023914317289, 18.215228855601, 0.2543411199572481, 0.2543411199572481),
    @("FAPs", "Fgf17", "Fgfr2", "FAPs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 3.344335848130888, 30.09902263317799, 0.4202757586428614, 0.4202757586428614),
    @("FAPs", "Fgf17", "Fgfr2", "MuSCs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.4494826044196666, 4.045343439777, 0.0564855478479598, 0.05648554784795982),
    @("MuSCs", "Fgf17", "Fgfr2", "ECs", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 0.597513709518, 5.377623385662, 0.07508831019693676, 0.07508831019693678),
    @("MuSCs", "Fgf17", "Fgfr2", "FAPs", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 0.9873375080262221, 8.886037572235999, 0.124076659482087, 0.124076659482087),
    @("MuSCs", "Fgf17", "Fgfr2", "MuSCs", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.1326993025526667, 1.194293722974, 0.0166760464810584, 0.0166760464810584),
    @("Resolving-Mac", "Fgf17", "Fgfr2", "ECs", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 1.414593, 4.243779, 0.3478871232761722, 0.3478871232761722, 0.146876719659, 1.321890476931, 0.01845769312198731, 0.01845769312198732),
    @("Resolving-Mac", "Fgf17", "Fgfr2", "FAPs", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 2.337487333333333, 7.012461999999999, 0.5748520910875596, 0.5748520910875596, 0.2427005306575555, 2.184304775918, 0.03049967296261124, 0.03049967296261125),
    @("Resolving-Mac", "Fgf17", "Fgfr2", "MuSCs", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 0.314161, 0.942483, 0.07726078563626818, 0.07726078563626819, 0.03261923190966667, 0.293573087187, 0.004099191307249969, 0.004099191307249971),
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $rowData[$j]
    }
}
